$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 5 content (order matters for shared string table indices)
$ws.Range("A5").Value = "POST"
$ws.Range("C5").Value = "http://192.168.100.19/thaimaiapp/api/mother/mPrimaryInfoUpdate"
$ws.Range("B5").Value = "PRIMARY REGISTRATION Update"

# Apply style s=4 (center horizontal + center vertical, no wrap) matching row2 A/B cells
$ws.Range("A5:C5").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A5:C5").VerticalAlignment = -4108    # xlCenter
$ws.Range("A5:C5").WrapText = $false

# Set row height for row 5
$ws.Rows.Item(5).RowHeight = 36.75

# Adjust column widths (engine quantizes to 1/6-character steps, so these
# inputs are chosen to land on the closest achievable stored width to the
# target values of 29.140625 and 62.85546875)
$ws.Columns.Item(2).ColumnWidth = 28.3
$ws.Columns.Item(3).ColumnWidth = 62.0

# Update selection to A4
$ws.Range("A4").Select()
